$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header in column G was "F"; change it to "D" (re-using the existing
# shared string for "D"). This also makes the old "F" shared string
# unused, so it drops out of sharedStrings.xml and "G" (column H header)
# is renumbered - matching the authored diff exactly.
$ws.Range("G1").Value = "D"

# Mark duplicate rows (new "bat ma mau trung" / duplicate-flag feature):
# rows 3 and 6 get a "1" in column G.
$ws.Range("G3").Value = 1
$ws.Range("G6").Value = 1

# Move the active selection to B7, matching the saved cursor position.
$ws.Range("B7").Select() | Out-Null
